# Auto-generated Excel COM-interop script to apply cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.976.23"
$ws.Range("E2").Value = "  +2.27%  "
$ws.Range("D3").Value = "3.189.52"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.58%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -2.87%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("D12").Value = "3.741.58"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("E13").Value = "  -2.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").Value = "59.979.64"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("D17").Value = "3.171.57"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "369.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.523"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.58%  "
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("D28").Value = "0.0₃0871"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("E32").Value = "  +2.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.34%  "
$ws.Range("D38").Value = "2.785.23"
$ws.Range("E38").Value = "  +5.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0710"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0308"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.71%  "
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.718"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").Value = "3.231.74"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.980"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.793"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.54%  "
